# Addressed QA for example_config file
#
# The last three rows of the "Data_Extent" table (years 2020-2022, rows
# 75-77) had their Year / Extent-on-land / Extent-at-sea / Total-extent
# values removed, leaving the cells blank but keeping their number
# formatting/styles. The sheet view was also scrolled down and the blank
# range left selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the values for rows 75-77 (columns A:D) while preserving the
# existing cell formatting (style indexes stay untouched).
$ws.Range("A75:D77").ClearContents()

# Leave the cleared range selected, as it was when the edit was made.
$ws.Range("A75:D77").Select()

# Scroll the window down so row 58 is at the top of the visible area.
$excel.ActiveWindow.ScrollRow = 58
$excel.ActiveWindow.ScrollColumn = 1
